# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51)
# to reflect the latest scrape from coinranking.com.
#
# A handful of "Price" (column D) values are purely numeric-looking
# strings (e.g. "245.05") that Excel would otherwise auto-convert to
# a floating point Number on assignment (losing the original fixed
# decimal-place text, e.g. "245.05" -> 245.0500000000001). Those are
# written with a leading apostrophe so Excel stores/keeps them as Text,
# matching the workbook's existing convention of storing every
# Price/Volume cell as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.867.42'
$ws.Range('E2').Value = '  +2.01%  '

# Row 3
$ws.Range('D3').Value = '1.893.95'
$ws.Range('E3').Value = '  +1.22%  '

# Row 4
$ws.Range('D4').Value = '''1.000'
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').Value = '''245.05'
$ws.Range('E5').Value = '  +4.29%  '

# Row 6
$ws.Range('D6').Value = '''0.9994'
$ws.Range('E6').Value = '  -0.10%  '

# Row 7
$ws.Range('D7').Value = '0.4778'
$ws.Range('E7').Value = '  +1.68%  '

# Row 8
$ws.Range('D8').Value = '''0.2910'
$ws.Range('E8').Value = '  +2.00%  '

# Row 9
$ws.Range('D9').Value = '''42.90'
$ws.Range('E9').Value = '  +3.02%  '

# Row 10
$ws.Range('D10').Value = '''0.06575'
$ws.Range('E10').Value = '  +0.31%  '

# Row 11
$ws.Range('D11').Value = '21.52'
$ws.Range('E11').Value = '  +1.10%  '

# Row 12
$ws.Range('D12').Value = '''0.07787'
$ws.Range('E12').Value = '  -0.47%  '

# Row 13
$ws.Range('D13').Value = '''97.09'
$ws.Range('E13').Value = '  +0.35%  '

# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.888.29'
$ws.Range('E14').Value = '  +0.66%  '

# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '''0.7371'
$ws.Range('E15').Value = '  +6.54%  '

# Row 16
$ws.Range('D16').Value = '''5.195'
$ws.Range('E16').Value = '  +2.34%  '

# Row 17
$ws.Range('D17').Value = '282.56'
$ws.Range('E17').Value = '  +5.23%  '

# Row 18
$ws.Range('D18').Value = '30.856.17'
$ws.Range('E18').Value = '  +2.00%  '

# Row 19
$ws.Range('D19').Value = '''13.63'
$ws.Range('E19').Value = '  -1.07%  '

# Row 20
$ws.Range('D20').Value = '''0.000007647'
$ws.Range('E20').Value = '  -0.74%  '

# Row 21
$ws.Range('D21').Value = '''1.000'
$ws.Range('E21').Value = '  +0.00%  '

# Row 22
$ws.Range('D22').Value = '2.161.68'
$ws.Range('E22').Value = '  +2.58%  '

# Row 23
$ws.Range('D23').Value = '''5.316'
$ws.Range('E23').Value = '  +1.21%  '

# Row 24
$ws.Range('D24').Value = '''1.000'
$ws.Range('E24').Value = '  +0.04%  '

# Row 25
$ws.Range('D25').Value = '''6.264'
$ws.Range('E25').Value = '  +1.49%  '

# Row 26
$ws.Range('D26').Value = '''9.364'
$ws.Range('E26').Value = '  -1.39%  '

# Row 27
$ws.Range('D27').Value = '''165.86'
$ws.Range('E27').Value = '  -0.03%  '

# Row 28
$ws.Range('D28').Value = '''19.21'
$ws.Range('E28').Value = '  +2.15%  '

# Row 29
$ws.Range('D29').Value = '''1.997'
$ws.Range('E29').Value = '  +3.10%  '

# Row 30
$ws.Range('D30').Value = '1.385'
$ws.Range('E30').Value = '  +1.03%  '

# Row 31
$ws.Range('D31').Value = '0.1002'
$ws.Range('E31').Value = '  +0.94%  '

# Row 32
$ws.Range('D32').Value = '''1.519'
$ws.Range('E32').Value = '  +4.28%  '

# Row 33
$ws.Range('D33').Value = '''4.395'
$ws.Range('E33').Value = '  +0.79%  '

# Row 34
$ws.Range('D34').Value = '''4.140'
$ws.Range('E34').Value = '  +1.97%  '

# Row 35
$ws.Range('D35').Value = '''0.04784'
$ws.Range('E35').Value = '  +0.79%  '

# Row 36
$ws.Range('D36').Value = '''1.134'
$ws.Range('E36').Value = '  +0.10%  '

# Row 37
$ws.Range('D37').Value = '''0.7078'
$ws.Range('E37').Value = '  +0.92%  '

# Row 38
$ws.Range('D38').Value = '''2.717'
$ws.Range('E38').Value = '  -0.07%  '

# Row 39
$ws.Range('E39').Value = '  +0.56%  '

# Row 40
$ws.Range('D40').Value = '''2.770'
$ws.Range('E40').Value = '  -0.41%  '

# Row 41
$ws.Range('D41').Value = '''6.491'
$ws.Range('E41').Value = '  +3.32%  '

# Row 42
$ws.Range('D42').Value = '''71.10'
$ws.Range('E42').Value = '  -2.56%  '

# Row 43
$ws.Range('D43').Value = '''1.939'
$ws.Range('E43').Value = '  -0.06%  '

# Row 44
$ws.Range('D44').Value = '''0.4230'
$ws.Range('E44').Value = '  +1.56%  '

# Row 45
$ws.Range('D45').Value = '''0.8480'
$ws.Range('E45').Value = '  +1.59%  '

# Row 46
$ws.Range('D46').Value = '''0.9992'
$ws.Range('E46').Value = '  -0.13%  '

# Row 47
$ws.Range('D47').Value = '102.91'
$ws.Range('E47').Value = '  -0.26%  '

# Row 48
$ws.Range('D48').Value = '''9.461'
$ws.Range('E48').Value = '  +3.09%  '

# Row 49
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '''956.22'
$ws.Range('E49').Value = '  -2.75%  '

# Row 50
$ws.Range('B50').Value = 'Aptos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D50').Value = '''7.156'
$ws.Range('E50').Value = '  +0.73%  '

# Row 51
$ws.Range('D51').Value = '''35.37'
$ws.Range('E51').Value = '  +2.42%  '
